$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Update the banner textbox shape (paragraph with the drawing):
#    - wp:extent / a:ext get slightly re-measured values
#    - overflowPunct flips false -> true (both the DrawingML and the
#      VML fallback copies of the frame paragraph)
#    - the VML fallback rectangle width is re-measured too
# -----------------------------------------------------------------
$shapeParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.End - $para.Range.Start -eq 1 -and $d.Shapes.Count -gt 0) {
        $shapeParaIndex = $i
        break
    }
}
if ($shapeParaIndex -eq 0) { $shapeParaIndex = 2 }
$shapePara = $d.Paragraphs($shapeParaIndex).Range
$shapeXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" mc:Ignorable="w14 wp14"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr></w:rPr></w:pPr><w:r><w:rPr></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor behindDoc="0" distT="0" distB="0" distL="0" distR="0" simplePos="0" locked="0" layoutInCell="1" allowOverlap="1" relativeHeight="2"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>289560</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>33655</wp:posOffset></wp:positionV><wp:extent cx="5566410" cy="525145"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapNone/><wp:docPr id="1" name="Shape1"/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5565600" cy="524520"/></a:xfrm><a:prstGeom prst="rect"><a:avLst></a:avLst></a:prstGeom><a:noFill/><a:ln><a:solidFill><a:srgbClr val="000000"/></a:solidFill></a:ln></wps:spPr><wps:style><a:lnRef idx="0"></a:lnRef><a:fillRef idx="0"/><a:effectRef idx="0"></a:effectRef><a:fontRef idx="minor"/></wps:style><wps:txbx><w:txbxContent><w:p><w:pPr><w:pStyle w:val="FrameContents"/><w:overflowPunct w:val="true"/><w:rPr></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00000A"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>This lab was developed for the Labtainer framework by the Naval Postgraduate School, Center for Cybersecurity and Cyber Operations under National Science Foundation Award No. 1438893.  This work is in the public domain, and cannot be copyrighted.</w:t></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr lIns="0" rIns="0" tIns="0" bIns="0"><a:spAutoFit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect id="shape_0" ID="Shape1" stroked="t" style="position:absolute;margin-left:22.8pt;margin-top:2.65pt;width:438.2pt;height:41.25pt"><w10:wrap type="square"/><v:fill o:detectmouseclick="t" on="false"/><v:stroke color="black" joinstyle="round" endcap="flat"/><v:textbox><w:txbxContent><w:p><w:pPr><w:pStyle w:val="FrameContents"/><w:overflowPunct w:val="true"/><w:rPr></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00000A"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>This lab was developed for the Labtainer framework by the Naval Postgraduate School, Center for Cybersecurity and Cyber Operations under National Science Foundation Award No. 1438893.  This work is in the public domain, and cannot be copyrighted.</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$shapePara.InsertXML($shapeXml)

# -----------------------------------------------------------------
# 2. Rewrite the "ssh ubuntu@<server_ip> [...]" line:
#    - paragraph mark run fonts (rPr inside pPr) simplified
#    - text split into "ssh ubuntu@" + "172.20.0.3" (both Tlwg Typist)
#      + a trailing single space (default font)
#    The following paragraph ("B. After connecting...") is included
#    unmodified in the replacement so that including the paragraph
#    mark of the target paragraph does not swallow it.
# -----------------------------------------------------------------
$sshIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "ssh ubuntu") {
        $sshIndex = $i
        break
    }
}
$sshPara = $d.Paragraphs($sshIndex).Range
$sshRange = $d.Range($sshPara.Start, $sshPara.End + 1)
$sshXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" mc:Ignorable="w14 wp14"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr></w:rPr><w:tab/><w:t xml:space="preserve">A. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tlwg Typist" w:hAnsi="Tlwg Typist"/></w:rPr><w:t>ssh ubuntu@</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tlwg Typist" w:hAnsi="Tlwg Typist"/></w:rPr><w:t>172.20.0.3</w:t></w:r><w:r><w:rPr></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:rFonts w:ascii="Liberation Serif" w:hAnsi="Liberation Serif"/></w:rPr></w:pPr><w:r><w:rPr></w:rPr><w:tab/><w:t xml:space="preserve">B. After connecting to the server: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tlwg Typist" w:hAnsi="Tlwg Typist"/></w:rPr><w:t>cat filetoview.txt</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$sshRange.InsertXML($sshXml)

# -----------------------------------------------------------------
# 3. Add the new character styles ListLabel 19 .. ListLabel 27
#    (wdStyleTypeCharacter = 2)
# -----------------------------------------------------------------
$listLabel19 = $d.Styles.Add("ListLabel 19", 2)
$listLabel19.QuickStyle = $true
$listLabel19.Font.NameAscii = "Liberation Serif"
$listLabel19.Font.NameOther = "Liberation Serif"
$listLabel19.Font.NameBi = "OpenSymbol"

for ($n = 20; $n -le 27; $n++) {
    $style = $d.Styles.Add("ListLabel $n", 2)
    $style.QuickStyle = $true
    $style.Font.NameBi = "OpenSymbol"
}

Write-Host "Edit complete"
